$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $val)
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" "57.485.87"
Set-TextValue "E2" "  -0.38%  "
Set-TextValue "D3" "3.108.26"
Set-TextValue "E3" "  +0.17%  "
Set-TextValue "E4" "  -0.01%  "
Set-TextValue "D5" "526.25"
Set-TextValue "E5" "  +0.57%  "
Set-TextValue "D6" "137.03"
Set-TextValue "E6" "  -3.26%  "
Set-TextValue "D7" "1.00"
Set-TextValue "E7" "  +0.03%  "
Set-TextValue "D8" "3.106.03"
Set-TextValue "E8" "  +0.11%  "
Set-TextValue "D9" "0.448"
Set-TextValue "E9" "  +2.14%  "
Set-TextValue "D10" "7.29"
Set-TextValue "E10" "  +1.18%  "
Set-TextValue "E11" "  -0.59%  "
Set-TextValue "E12" "  +2.76%  "
Set-TextValue "D13" "3.644.54"
Set-TextValue "E13" "  +0.19%  "
Set-TextValue "E14" "  +2.63%  "
Set-TextValue "D15" "25.35"
Set-TextValue "E15" "  -1.59%  "
Set-TextValue "E16" "  -0.47%  "
Set-TextValue "D17" "57.595.15"
Set-TextValue "E17" "  -0.35%  "
Set-TextValue "D18" "3.108.44"
Set-TextValue "E18" "  +0.00%  "
Set-TextValue "E19" "  -2.50%  "
Set-TextValue "D20" "12.40"
Set-TextValue "E20" "  -2.86%  "
Set-TextValue "E21" "  -2.08%  "
Set-TextValue "D22" "345.53"
Set-TextValue "E22" "  +2.30%  "
Set-TextValue "E23" "  -0.05%  "
Set-TextValue "D24" "67.70"
Set-TextValue "E24" "  +1.77%  "
Set-TextValue "E25" "  -2.10%  "
Set-TextValue "E26" "  -1.38%  "
Set-TextValue "D27" "1.00"
Set-TextValue "E27" "  -0.03%  "
Set-TextValue "D28" "0.0₃0894"
Set-TextValue "E28" "  -2.23%  "
Set-TextValue "D29" "7.42"
Set-TextValue "E29" "  +3.62%  "
Set-TextValue "E30" "  +0.04%  "
Set-TextValue "E31" "  +0.24%  "
Set-TextValue "E32" "  -6.76%  "
Set-TextValue "D33" "20.82"
Set-TextValue "E33" "  -0.38%  "
Set-TextValue "D34" "4.95"
Set-TextValue "E34" "  +7.64%  "
Set-TextValue "E35" "  -3.12%  "
Set-TextValue "D36" "158.46"
Set-TextValue "E36" "  +1.62%  "
Set-TextValue "D37" "6.07"
Set-TextValue "E37" "  -0.70%  "
Set-TextValue "D38" "25.88"
Set-TextValue "E38" "  -4.18%  "
Set-TextValue "E39" "  -1.51%  "
Set-TextValue "E40" "  +5.83%  "
Set-TextValue "E41" "  +0.05%  "
Set-TextValue "E42" "  +3.78%  "
Set-TextValue "D43" "0.700"
Set-TextValue "E43" "  +2.63%  "
Set-TextValue "D44" "3.149.01"
Set-TextValue "E44" "  +0.19%  "
Set-TextValue "D45" "2.378.81"
Set-TextValue "E45" "  +3.56%  "
Set-TextValue "D46" "36.60"
Set-TextValue "E46" "  -0.43%  "
Set-TextValue "E47" "  +0.02%  "
Set-TextValue "E48" "  +3.33%  "
Set-TextValue "D49" "0.975"
Set-TextValue "E49" "  -0.82%  "
Set-TextValue "E50" "  -0.52%  "
Set-TextValue "D51" "19.80"
Set-TextValue "E51" "  -3.42%  "
